# Switching to Summer time
# Shift every timestamp in column A forward by 3 days (the data window moved
# from 2025-03-28 to 2025-03-31) and update the corresponding Notified
# Production (MW) values in column B for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 97

# Shift all timestamps (column A) forward by 3 days, keeping the existing
# number formatting/style on each cell untouched.
for ($r = 2; $r -le $lastRow; $r++) {
    $oldDate = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $oldDate + 3
}

# New Notified Production (MW) readings for the shifted day (rows 22-77).
$bValues = @{
    22 = 40
    23 = 31
    24 = 29
    25 = 32
    26 = 134
    27 = 140
    28 = 149
    29 = 157
    30 = 302
    31 = 314
    32 = 327
    33 = 343
    34 = 476
    35 = 488
    36 = 498
    37 = 509
    38 = 589
    39 = 597
    40 = 605
    41 = 611
    42 = 655
    43 = 658
    44 = 661
    45 = 664
    46 = 658
    47 = 658
    48 = 656
    49 = 651
    50 = 609
    51 = 601
    52 = 593
    53 = 582
    54 = 521
    55 = 512
    56 = 502
    57 = 493
    58 = 393
    59 = 384
    60 = 373
    61 = 364
    62 = 262
    63 = 252
    64 = 243
    65 = 235
    66 = 116
    67 = 110
    68 = 103
    69 = 99
    70 = 10
    71 = 8
    72 = 7
    73 = 6
    74 = 1
    75 = 0
    76 = 0
    77 = 0
}

foreach ($r in $bValues.Keys) {
    $ws.Cells.Item($r, 2).Value = $bValues[$r]
}
